# Fixed wrong setting of age in dismantling module
$wb = $excel.ActiveWorkbook

# --- Sheet "times": update StartTime / StopTime ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 45657.99861111111
$wsTimes.Range("B3").Value = 46021.99861111111

# --- Sheet "scenario_data_emlab": shrink from 11 year-columns (2027-2037) to 2 (2024-2025) ---
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")

# Drop the now-unused year columns D:L (columns 2030-2037 onward) entirely
$wsScenario.Range("D1:L8").Clear()

# Update the remaining year headers and values
$wsScenario.Range("B1").Value = 2024
$wsScenario.Range("C1").Value = 2025

$wsScenario.Range("B2").Value = 40.3
$wsScenario.Range("B5").Value = 14.32
$wsScenario.Range("B6").Value = 21.7
$wsScenario.Range("B7").Value = 54.81

# --- Sheet "renewables": drop the placeholder OtherPV rows (20222100086 etc.) and keep
#     only the real RunOfRiver / OtherPV / WindOff plants, moved up into rows 17-19 ---
$wsRenew = $wb.Worksheets.Item("renewables")

# Remove rows 20-39 entirely (shifts nothing below them, they are the last rows)
$wsRenew.Range("A20:A39").EntireRow.Delete()

# Overwrite rows 17-19 with the data that used to live in rows 37-39
$wsRenew.Range("B17").Value = 20151200026
$wsRenew.Range("C17").Value = 8858.749999999998
$wsRenew.Range("E17").Value = "RunOfRiver"

$wsRenew.Range("B18").Value = 20152100030
$wsRenew.Range("C18").Value = 53555.51607579708

$wsRenew.Range("B19").Value = 20152300031
$wsRenew.Range("C19").Value = 10271.8
$wsRenew.Range("D19").Value = 2.7
$wsRenew.Range("E19").Value = "WindOff"
